# fix: unique command names in XLSX - prefix protocol name to each step
# For every "protocol" worksheet (named after the scenario it represents),
# prefix the worksheet's own name + a space onto the existing value of each
# data row's column A cell (the "Step..." / label names), leaving the header
# row (row 1, "Name") and all other columns untouched.

$wb = $excel.ActiveWorkbook

$sheetNames = @(
    "price1", "price2",
    "discount1", "discount2",
    "free1", "free2",
    "nomoney1", "nomoney2",
    "noppv1", "noppv2",
    "card1", "card2",
    "nosex1", "nosex2",
    "offtopic1", "offtopic2",
    "real1", "real2",
    "voice1", "voice2",
    "customyes1", "customyes2",
    "customno1", "customno2",
    "done1", "done2",
    "cumcontrol", "dickpic", "boosters"
)

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $lastRow = $ws.UsedRange.Rows.Count()

    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 1)
        $current = $cell.Value()
        if ($current -ne $null -and $current -ne "") {
            $cell.Value = $sheetName + " " + $current
        }
    }
}
